# Actualización automática de catálogo y fotos
# The "San Francisco Velvet Topo" product (row 24) is discontinued and
# removed from the catalogo sheet. Deleting the whole row shifts every
# row below it up by one, which is exactly what the target workbook shows
# (tables, data validations and sharedStrings all renumber naturally).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

$ws.Rows.Item(24).Delete()
